$wb = $excel.ActiveWorkbook

# Rename "Bottrell pooled" -> "bottrell_pooled"
$bottrellPooled = $wb.Worksheets.Item("Bottrell pooled")
$bottrellPooled.Name = "bottrell_pooled"

# Make it the active sheet/tab and move the selection to H28
# (this also flips the previously-active "McCauley rotifers" tab back to unselected)
$bottrellPooled.Activate()
$bottrellPooled.Range("H28").Select()
